$wb = $excel.ActiveWorkbook

# --- Sheet "TestCase": insert a new "Note" column before current column G (TesterName) ---
$ws1 = $wb.Worksheets.Item("TestCase")

# Insert a new column at G; everything from G onward (TesterName..UpdatedDateTime) shifts
# right by one. The new cell G1 automatically inherits the neighboring header's style (s=4),
# and the existing dataValidation / named ranges referencing columns to the right shift too.
$ws1.Columns("G").Insert()

# Header text for the newly inserted column
$ws1.Range("G1").Value = "Note"

# Give the new column an explicit (non-bestFit) width, matching the target layout
$ws1.Columns("G").ColumnWidth = 11

# The existing AutoFilter range doesn't auto-grow with the inserted column, so drop it and
# reapply over the full, now-one-column-wider header row.
$ws1.AutoFilterMode = $false
$ws1.Range("A1:R1").AutoFilter()

# The workbook-level hidden _FilterDatabase name still points at the old A1:Q1 range -
# repoint it at the new A1:R1 range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "TestCase!_FilterDatabase") {
        $n.RefersTo = "=TestCase!`$A`$1:`$R`$1"
    }
}

# --- Sheet "Note": bump template version / metadata, mention the renamed column ---
$ws2 = $wb.Worksheets.Item("Note")

# Widen column C and stop relying on bestFit for its width
$ws2.Columns("C").ColumnWidth = 39.833333333333336

# Updated datetime (serial date)
$ws2.Range("B3").Value = 43753

# Total columns count
$ws2.Range("B4").Value = 18

# Template version bump
$ws2.Range("B2").Value = "v1.3.1.0"

# Notes text body - leading apostrophe preserves the quote-prefixed text formatting that the
# original cell already had (it starts with "-", which Excel treats as needing a quote prefix).
$ws2.Range("C4").Formula = "'- Do not change ordinal of columns in the left of column System Validation. `n- Can change name of any columns.`n- Can add or remove columns in the right of column System Validation."
